$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).Style = "Normal"
}

# Rows with updated Price (D) and/or Volume(1h) (E) values
Set-TextValue $ws "D2" "56.760.10"
$ws.Range("E2").Value = "  +3.33%  "
Set-TextValue $ws "D3" "2.326.00"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws "D5" "521.39"
$ws.Range("E5").Value = "  +3.14%  "
Set-TextValue $ws "D6" "135.31"
$ws.Range("E6").Value = "  +4.46%  "
Set-TextValue $ws "D7" "0.995"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.39%  "
Set-TextValue $ws "D9" "2.350.38"
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("E11").Value = "  -0.94%  "
Set-TextValue $ws "D12" "5.31"
$ws.Range("E12").Value = "  +2.96%  "
Set-TextValue $ws "D13" "0.344"
$ws.Range("E13").Value = "  +0.88%  "
Set-TextValue $ws "D14" "24.08"
$ws.Range("E14").Value = "  +1.64%  "
Set-TextValue $ws "D15" "2.741.24"
$ws.Range("E15").Value = "  +1.52%  "
Set-TextValue $ws "D16" "56.810.35"
$ws.Range("E16").Value = "  +3.42%  "
Set-TextValue $ws "D17" "0.0000135"
$ws.Range("E17").Value = "  +2.22%  "
Set-TextValue $ws "D18" "2.326.84"
$ws.Range("E18").Value = "  +1.45%  "
Set-TextValue $ws "D19" "10.54"
$ws.Range("E19").Value = "  +0.38%  "
Set-TextValue $ws "D20" "4.24"
$ws.Range("E20").Value = "  +1.63%  "
Set-TextValue $ws "D21" "323.69"
$ws.Range("E21").Value = "  +3.83%  "
Set-TextValue $ws "D22" "6.62"
$ws.Range("E22").Value = "  -0.46%  "
Set-TextValue $ws "D23" "0.999"
$ws.Range("E23").Value = "  +0.18%  "
Set-TextValue $ws "D24" "60.67"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  +9.55%  "
Set-TextValue $ws "D26" "0.991"
$ws.Range("E26").Value = "  -0.18%  "
Set-TextValue $ws "D27" "8.02"
$ws.Range("E27").Value = "  +6.84%  "
Set-TextValue $ws "D28" "1.30"
$ws.Range("E28").Value = "  +14.11%  "
Set-TextValue $ws "D29" "0.0₃0747"
$ws.Range("E29").Value = "  +5.65%  "
$ws.Range("E30").Value = "  +5.25%  "
Set-TextValue $ws "D31" "166.99"
$ws.Range("E31").Value = "  -2.87%  "
Set-TextValue $ws "D32" "6.23"
$ws.Range("E32").Value = "  +1.25%  "
Set-TextValue $ws "D33" "18.40"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  +2.61%  "
Set-TextValue $ws "D37" "0.932"
$ws.Range("E37").Value = "  +0.65%  "
Set-TextValue $ws "D38" "4.06"
$ws.Range("E38").Value = "  +4.80%  "
Set-TextValue $ws "D39" "1.57"
$ws.Range("E39").Value = "  +8.04%  "
Set-TextValue $ws "D40" "37.95"
$ws.Range("E40").Value = "  +2.99%  "
Set-TextValue $ws "D41" "0.380"
$ws.Range("E41").Value = "  +0.90%  "
Set-TextValue $ws "D44" "5.32"
$ws.Range("E44").Value = "  +3.92%  "
Set-TextValue $ws "D45" "280.08"
$ws.Range("E45").Value = "  +7.24%  "
$ws.Range("E46").Value = "  +2.22%  "
Set-TextValue $ws "D47" "0.0507"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E51").Value = "  +0.96%  "

# Row 42 <-> Row 43 swap (Aave and Filecoin swapped ranking positions, with updated prices)
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D42" "3.62"
$ws.Range("E42").Value = "  +5.30%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D43" "139.01"
$ws.Range("E43").Value = "  +3.61%  "

# Row 49 <-> Row 50 swap (VeChain and InjectiveProtocol swapped ranking positions, with updated prices)
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D49" "18.18"
$ws.Range("E49").Value = "  +9.94%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D50" "0.0219"
$ws.Range("E50").Value = "  +3.40%  "
